# Update analysis of false positives and negatives
# - Rename "Comment" header (E1) to "Comment LD"
# - Add a new "Comment SB" column (F) with per-row remarks
# - Re-apply the AutoFilter over the header row only (A1:F1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing comment column header
$ws.Range("E1").Value = "Comment LD"

# New column header
$ws.Range("F1").Value = "Comment SB"

# New column data (rows 2-12)
$ws.Range("F2").Value = "Matching due to BELOW strategy before ABOVE strategy in combination with matched text block containing link. Can't be fixed without breaking other test cases."
$ws.Range("F3").Value = "Syntactically not similar enough"
$ws.Range("F4").Value = "Syntactically not similar enough"
$ws.Range("F5").Value = "Fixed (now using edit-based metric if string only contain one token)"
$ws.Range("F6").Value = "Matching strategy fails here, because semantic and syntactical similarity differ"
$ws.Range("F7").Value = "Matching strategy fails here, because semantic and syntactical similarity differ"
$ws.Range("F8").Value = "Matching strategy fails here, because semantic and syntactical similarity differ"
$ws.Range("F9").Value = "false positive"
$ws.Range("F10").Value = "Matching strategy fails here, because semantic and syntactical similarity differ"
$ws.Range("F11").Value = "Matching strategy fails here, because semantic and syntactical similarity differ"
$ws.Range("F12").Value = "Matching strategy fails here, because semantic and syntactical similarity differ"

# Mirror the highlighted formatting on the new column for the rows that were
# already highlighted in the "Comment LD" column (2, 5, 9)
$ws.Range("F2").Style = "Normal 2"
$ws.Range("F2").Interior.ColorIndex = -4142
$ws.Range("F5").Style = "Normal 2"
$ws.Range("F5").Interior.ColorIndex = -4142
$ws.Range("F9").Style = "Normal 2"
$ws.Range("F9").Interior.ColorIndex = -4142

# Column width for the new column (approximate best-fit width)
$ws.Columns.Item(6).ColumnWidth = 136.5

# Re-point the AutoFilter to the header row only (A1:F1) instead of the
# full previous data range (A1:E84)
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$ws.Range("A1:F1").AutoFilter()

# The hidden _FilterDatabase defined name needs to follow the new
# AutoFilter range as well.
$fd = $wb.Names.Item(1)
$fd.RefersTo = "=Sheet1!`$A`$1:`$F`$1"

# Selection moves to A12
$ws.Range("A12").Select()
